$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the H1 title.
#    Pattern: empty run, bold "Meta description" run, then a plain run with
#    the rest of the sentence (matches the style used elsewhere in the doc).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
$titlePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$insertionPoint = $newPara.Range
$insertionPoint.Collapse(1)

$metaXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
    "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:body>" +
    "<w:p>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
    "<w:r><w:t>: Get to know the features and functions of Divine Showdown, a slot game from Play 'N Go. Play for free and read our review.</w:t></w:r>" +
    "</w:p>" +
    "</w:body>" +
    "</w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

$insertionPoint.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph that used to sit near the
#    bottom of the document (right before the italic meta-description blurb).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "Play Divine Showdown Free | Review of Play 'N Go Slot Game") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -ge 1) {
    $d.Paragraphs.Item($targetIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Replace the text of the trailing italic paragraph (formerly the meta
#    description blurb) with the new DALLE image prompt, preserving the
#    leading empty run and the italic run formatting.
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count2)
$targetRange = $lastPara.Range
$targetRange.MoveEnd(1, -1)

$newPrompt = "Prompt for DALLE: Create a cartoon-style feature image for the online slot game `"Divine Showdown`" that features a happy Maya warrior wearing glasses. The image should capture the divine powers and epic background of the game, incorporating the four deities in an exciting and engaging way. Use bright, bold colors to capture the attention of potential players and showcase the game's excitement. Be sure to include the game's title and the Play 'N Go logo to effectively promote the game."

$targetRange.Text = $newPrompt

Write-Output "Edit complete."
